# This workbook (weekly price listing for "Perejil" at Vega Central
# Mapocho de Santiago) gets a new weekly observation inserted as row 351,
# pushing every existing row from 351..428 down by one (to 352..429).
#
# Insert a new row at 351 (this shifts rows 351..428 -> 352..429 and
# grows the used range to A1:R429 automatically).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(351).Insert()

# Populate the newly inserted row 351 with the new observation.
$ws.Cells.Item(351, 1).Value()  = 9
$ws.Cells.Item(351, 2).Value()  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(351, 3).Value()  = "Metropolitana"
$ws.Cells.Item(351, 4).Value()  = 44798
$ws.Cells.Item(351, 5).Value()  = 13
$ws.Cells.Item(351, 6).Value()  = 100112044
$ws.Cells.Item(351, 7).Value()  = "Perejil"
$ws.Cells.Item(351, 8).Value()  = "Sin especificar"
$ws.Cells.Item(351, 9).Value()  = "Primera"
$ws.Cells.Item(351, 10).Value() = 70
$ws.Cells.Item(351, 11).Value() = 12000
$ws.Cells.Item(351, 12).Value() = 13000
$ws.Cells.Item(351, 13).Value() = 12500
$ws.Cells.Item(351, 14).Value() = "`$/docena de atados"
$ws.Cells.Item(351, 15).Value() = "Región Metropolitana"
$ws.Cells.Item(351, 16).Value() = 4167
$ws.Cells.Item(351, 17).Value() = 3
$ws.Cells.Item(351, 18).Value() = "Hortaliza"
